# Applies the Case_1_243 (380 kV) vm_pu.xlsx rerun: bus B voltage setpoint
# moved from 1.05 to 1.02 p.u., and every other bus voltage in rows 2-25
# (columns B-F, I-N) was recomputed by the load-flow solver accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 1.02
    "C2" = 1.012646048344971
    "D2" = 1.049514405841244
    "E2" = 1.01466339299724
    "F2" = 1.051573407879197
    "I2" = 1.038803366697091
    "J2" = 1.01788756760643
    "K2" = 1.0522712518193
    "L2" = 1.017520870022129
    "M2" = 1.054324534896809
    "N2" = 1.010073930159328
    "B3" = 1.02
    "C3" = 1.013549519963768
    "D3" = 1.050089998678811
    "E3" = 1.015428324143953
    "F3" = 1.052314638827817
    "I3" = 1.038899930785042
    "J3" = 1.018425072773582
    "K3" = 1.052659717375361
    "L3" = 1.018090907799843
    "M3" = 1.054878622829201
    "N3" = 1.010255078974674
    "B4" = 1.02
    "C4" = 1.014134964783015
    "D4" = 1.050462293400085
    "E4" = 1.01592439557604
    "F4" = 1.052794532978703
    "I4" = 1.038961038847567
    "J4" = 1.018773122249264
    "K4" = 1.052910226330724
    "L4" = 1.018460212566834
    "M4" = 1.05523675501084
    "N4" = 1.010372282890692
    "B5" = 1.02
    "C5" = 1.014381285224642
    "D5" = 1.050618767394655
    "E5" = 1.016133208115818
    "F5" = 1.052996342327701
    "I5" = 1.038986398566641
    "J5" = 1.018919500447361
    "K5" = 1.05301533404406
    "L5" = 1.018615575347287
    "M5" = 1.055387216078205
    "N5" = 1.010421552222402
    "B6" = 1.02
    "C6" = 1.014422655169865
    "D6" = 1.050645037731904
    "E6" = 1.016168284101487
    "F6" = 1.053030230563736
    "I6" = 1.038990637192582
    "J6" = 1.018944081371784
    "K6" = 1.053032969946069
    "L6" = 1.018641667684639
    "M6" = 1.055412473366001
    "N6" = 1.010429824559692
    "B7" = 1.02
    "C7" = 1.014138255346696
    "D7" = 1.050464384368182
    "E7" = 1.015927184703638
    "F7" = 1.052797229325055
    "I7" = 1.038961379003549
    "J7" = 1.018775077934608
    "K7" = 1.052911631597771
    "L7" = 1.018462288112761
    "M7" = 1.055238765865416
    "N7" = 1.010372941243026
    "B8" = 1.02
    "C8" = 1.012951206490113
    "D8" = 1.049708960334244
    "E8" = 1.014921674308181
    "F8" = 1.051823853412143
    "I8" = 1.038836285251459
    "J8" = 1.01806916804702
    "K8" = 1.0524027113997
    "L8" = 1.017713422601904
    "M8" = 1.054511873081068
    "N8" = 1.010135152358495
    "B9" = 1.02
    "C9" = 1.01086595001112
    "D9" = 1.048376738660991
    "E9" = 1.013158406652278
    "F9" = 1.050110798263228
    "I9" = 1.038605357832151
    "J9" = 1.016827210280421
    "K9" = 1.051499463574233
    "L9" = 1.016397342479376
    "M9" = 1.053228021616085
    "N9" = 1.009716069758664
    "B10" = 1.02
    "C10" = 1.009480206423023
    "D10" = 1.047488006061471
    "E10" = 1.011988751501408
    "F10" = 1.048970352184437
    "I10" = 1.038444399525735
    "J10" = 1.016000611978144
    "K10" = 1.050893055294736
    "L10" = 1.015522391729899
    "M10" = 1.052370240413999
    "N10" = 1.009436662055396
    "B11" = 1.02
    "C11" = 1.008881229647502
    "D11" = 1.047103064505667
    "E11" = 1.011483686956871
    "F11" = 1.048476933331566
    "I11" = 1.038373051525875
    "J11" = 1.015643025308107
    "K11" = 1.050629489363753
    "L11" = 1.015144121257898
    "M11" = 1.051998390849902
    "N11" = 1.009315676442769
    "B12" = 1.02
    "C12" = 1.008658903373573
    "D12" = 1.046960065156247
    "E12" = 1.011296296024092
    "F12" = 1.048293718366304
    "I12" = 1.038346302348901
    "J12" = 1.01551025325309
    "K12" = 1.050531442699739
    "L12" = 1.015003704460805
    "M12" = 1.051860207593467
    "N12" = 1.009270737411528
    "B13" = 1.02
    "C13" = 1.008706585837371
    "D13" = 1.046990739640538
    "E13" = 1.011336482372077
    "F13" = 1.04833301571957
    "I13" = 1.038352051317941
    "J13" = 1.015538730958941
    "K13" = 1.050552480645043
    "L13" = 1.015033820261919
    "M13" = 1.051889851133868
    "N13" = 1.009280376960445
    "B14" = 1.02
    "C14" = 1.00886284880979
    "D14" = 1.047091244425948
    "E14" = 1.011468192819089
    "F14" = 1.048461787431135
    "I14" = 1.038370845471637
    "J14" = 1.015632049270542
    "K14" = 1.050621387775144
    "L14" = 1.01513251252481
    "M14" = 1.051986969827245
    "N14" = 1.009311961757433
    "B15" = 1.02
    "C15" = 1.00895914887148
    "D15" = 1.047153166832926
    "E15" = 1.011549372182659
    "F15" = 1.048541136321339
    "I15" = 1.038382392418432
    "J15" = 1.015689552626729
    "K15" = 1.050663824366693
    "L15" = 1.01519333200744
    "M15" = 1.052046799734717
    "N15" = 1.009331422261281
    "B16" = 1.02
    "C16" = 1.009519980961683
    "D16" = 1.047513551147912
    "E16" = 1.012022300687189
    "F16" = 1.049003107439357
    "I16" = 1.038449099923684
    "J16" = 1.01602435099742
    "K16" = 1.050910526650357
    "L16" = 1.015547508807783
    "M16" = 1.052394910071232
    "N16" = 1.009444691505307
    "B17" = 1.02
    "C17" = 1.009872060587145
    "D17" = 1.047739581673024
    "E17" = 1.012319333234354
    "F17" = 1.049292998990175
    "I17" = 1.038490502087838
    "J17" = 1.016234451861336
    "K17" = 1.051065013436487
    "L17" = 1.015769833098184
    "M17" = 1.052613158231265
    "N17" = 1.009515742562786
    "B18" = 1.02
    "C18" = 1.01007752480528
    "D18" = 1.047871410256304
    "E18" = 1.012492722637071
    "F18" = 1.049462126437051
    "I18" = 1.038514491820288
    "J18" = 1.016357032471509
    "K18" = 1.051155027601911
    "L18" = 1.015899567841708
    "M18" = 1.052740417583015
    "N18" = 1.009557185390018
    "B19" = 1.02
    "C19" = 1.010147600080009
    "D19" = 1.047916358427537
    "E19" = 1.0125518668501
    "F19" = 1.049519801001709
    "I19" = 1.038522644628141
    "J19" = 1.016398834742303
    "K19" = 1.051185703877111
    "L19" = 1.015943813624924
    "M19" = 1.05278380269903
    "N19" = 1.009571316300615
    "B20" = 1.02
    "C20" = 1.009834275202882
    "D20" = 1.047715331878579
    "E20" = 1.012287450442286
    "F20" = 1.049261892359374
    "I20" = 1.038486076513988
    "J20" = 1.016211906674948
    "K20" = 1.051048448317499
    "L20" = 1.01574597392759
    "M20" = 1.052589746506485
    "N20" = 1.009508119458941
    "B21" = 1.02
    "C21" = 1.008816828822374
    "D21" = 1.047061648665127
    "E21" = 1.011429401495493
    "F21" = 1.048423865619759
    "I21" = 1.038365317883223
    "J21" = 1.015604567921654
    "K21" = 1.05060110037392
    "L21" = 1.01510344762737
    "M21" = 1.051958372481621
    "N21" = 1.00930266080698
    "B22" = 1.02
    "C22" = 1.008178048358253
    "D22" = 1.046650566908443
    "E22" = 1.010891142481755
    "F22" = 1.047897329275079
    "I22" = 1.038287961445338
    "J22" = 1.015223008772851
    "K22" = 1.050318987989736
    "L22" = 1.014699985410084
    "M22" = 1.051561046180074
    "N22" = 1.009173483403991
    "B23" = 1.02
    "C23" = 1.008516589553684
    "D23" = 1.046868496439214
    "E23" = 1.011176366601103
    "F23" = 1.048176420821071
    "I23" = 1.038329104889876
    "J23" = 1.015425251804179
    "K23" = 1.050468620797968
    "L23" = 1.01491381871122
    "M23" = 1.051771709476937
    "N23" = 1.009241962393054
    "B24" = 1.02
    "C24" = 1.009851348476135
    "D24" = 1.047726289351752
    "E24" = 1.012301856485653
    "F24" = 1.049275947987857
    "I24" = 1.038488076733085
    "J24" = 1.016222093774457
    "K24" = 1.051055933676601
    "L24" = 1.015756754684581
    "M24" = 1.052600325384423
    "N24" = 1.009511564011651
    "B25" = 1.02
    "C25" = 1.011404263861287
    "D25" = 1.048721262702417
    "E25" = 1.013613228261023
    "F25" = 1.050553393945252
    "I25" = 1.038666296766876
    "J25" = 1.017148048830489
    "K25" = 1.051733730245339
    "L25" = 1.016737156624945
    "M25" = 1.053560268090373
    "N25" = 1.009824418002891
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
